$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.881.41'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '2.582.38'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.85'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.58'
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '2.588.82'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.55'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("E13").Value = '  +3.30%  '
$ws.Range("D14").Value = '3.035.20'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").Value = '58.939.79'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.36'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '2.581.00'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '336.79'
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.07'
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.49'
$ws.Range("E22").Value = '  +2.99%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.04'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.99'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '0.0₃0720'
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.92'
$ws.Range("E31").Value = '  -5.12%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.63'
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.02'
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.95'
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.11'
$ws.Range("E36").Value = '  -2.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.70'
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.820'
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("E40").Value = '  -7.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.49'
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '269.53'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.77'
$ws.Range("E44").Value = '  +1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0952'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.586'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0515'
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.30'
$ws.Range("E48").Value = '  -2.08%  '
$ws.Range("D49").Value = '1.959.79'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("E51").Value = '  -0.27%  '
